$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for rows 2-10, columns B-E and G (F unchanged)
$data = @{
    2  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248 }
    3  = @{ B = 0.1169995834814548; C = 0.3048912486333797;   D = 0.1496068669990043;  E = 0.5333859586016987;  G = 1.104883657715537 }
    4  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.7210945179870265;  E = 13.86384647080068;   G = 19.48425592650926 }
    5  = @{ B = 0.1169995834814548; C = 0.3048912486333797;   D = 0.7210945179870265;  E = 13.86384647080068;   G = 15.00683182090255 }
    6  = @{ B = 1.445647641019636;  C = 1.626987699542094;   D = 3.223369029078222;   E = 13.86384647080068;   G = 20.15985084044064 }
    7  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 18.71679738969934;   E = 0.5333859586016987;  G = 24.14949828602258 }
    8  = @{ B = 3.272327238179451;  C = 1.626987699542094;   D = 0.1496068669990043;  E = 0.5333859586016987;  G = 5.582307763322248 }
    9  = @{ B = 0.2881169905109251; C = 109.9114832445916;    D = 189.6080260415259;   E = 2797.565817734744;   G = 3097.373444011372 }
    10 = @{ B = 0.01253208636536152; C = 0.3048912486333797;  D = 18.71679738969934;   E = 2797.565817734744;   G = 2816.600038459442 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
